# "chore: update News.xlsx with new Moore-Czech releases"
# Append the latest Moore Czech press-release rows beneath the existing
# "News" table header on Sheet1 (table itself is left at its original C3:I4
# extent; the new rows simply land on the sheet below/around it).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newsItems = @(
    @{ Title = 'Online výuka bez kompromisů. Online Střední škola Spektrum potvrdila svou kvalitu maturitními výsledky ​'; PostDate = '30.06.2025'; Url = 'https://www.moore-czech.cz/tiskove-zpravy/cerven-2025/online-stredni-skola'; Source = 'Moore Czech s.r.o.' }
    @{ Title = 'Vysoké nároky, celoživotní vzdělávání i zapojení AI. Podoba moderního účetnictví se prudce mění'; PostDate = '11.06.2025'; Url = 'https://www.moore-czech.cz/tiskove-zpravy/cerven-2025/podoba-moderniho-ucetnictvi'; Source = 'Moore Czech s.r.o.' }
    @{ Title = 'Digitální škola patří mezi nejlepší v regionu. Umožňuje studium i profesionálním sportovcům'; PostDate = '27.05.2025'; Url = 'https://www.moore-czech.cz/tiskove-zpravy/kveten-2025/digitalni-skola-patri-mezi-nejlepsi-v-regionu'; Source = 'Moore Czech s.r.o.' }
    @{ Title = 'Energetika stále čelí řadě nejistot. Přechod na obnovitelné zdroje brzdí poplatky i infrastruktura'; PostDate = '16.04.2025'; Url = 'https://www.moore-czech.cz/tiskove-zpravy/duben-2025/energetika-stale-celi-rade-nejistot'; Source = 'Moore Czech s.r.o.' }
    @{ Title = 'Vztah mezi klientem a finanční institucí posilují digitální procesy i digitální peníze'; PostDate = '03.04.2025'; Url = 'https://www.moore-czech.cz/tiskove-zpravy/duben-2025/vztah-mezi-klientem-a-financni-instituci'; Source = 'Moore Czech s.r.o.' }
    @{ Title = 'Talk2amy vstupuje do nové éry díky strategickému partnerství s Moore Czech Republic'; PostDate = '02.04.2025'; Url = 'https://www.moore-czech.cz/tiskove-zpravy/duben-2025/talk2amy-vstupuje-do-nove-ery'; Source = 'Moore Czech s.r.o.' }
    @{ Title = 'Firmy v roce 2025 čekají další daňové změny, stavebnictví nebude výjimkou'; PostDate = '13.03.2025'; Url = 'https://www.moore-czech.cz/tiskove-zpravy/brezen-2025/danove-zmeny-2025'; Source = 'Moore Czech s.r.o.' }
    @{ Title = 'Home-office proměnil trh práce. Zůstane i nadále běžnou praxí?'; PostDate = '26.02.2025'; Url = 'https://www.moore-czech.cz/tiskove-zpravy/unor-2025/home-office-promenil-trh-prace-zustane-i-nadale-be'; Source = 'Moore Czech s.r.o.' }
    @{ Title = 'Skupina ANACOT CAPITAL opět rozšiřuje své působení. Do svého portfolia získala blanenského výrobce pistolí Laugo Arms Czechoslovakia'; PostDate = '09.01.2025'; Url = 'https://www.moore-czech.cz/tiskove-zpravy/leden-2025/skupina-anacot-capital-opet-rozsiruje'; Source = 'Moore Czech s.r.o.' }
    @{ Title = 'Skupina Moore Czech Republic expanduje na Slovensko v oblasti poradenských služeb'; PostDate = '19.12.2024'; Url = 'https://www.moore-czech.cz/tiskove-zpravy/prosinec-2024/skupina-moore-czech-republic-expanduje'; Source = 'Moore Czech s.r.o.' }
    @{ Title = 'Zůstat efektivní a nezbláznit se. Moore Czech Republic se připojuje ke světové iniciativě usilující o snižování stresu na pracovišti'; PostDate = '17.12.2024'; Url = 'https://www.moore-czech.cz/tiskove-zpravy/prosinec-2024/zustat-efektivni-a-nezblaznit-se'; Source = 'Moore Czech s.r.o.' }
    @{ Title = 'Rekvalifikace jako výtah do vyšších pater profesních i mzdových. Zájem roste, podpora nikoli'; PostDate = '05.12.2024'; Url = 'https://www.moore-czech.cz/tiskove-zpravy/prosinec-2024/rekvalifikace-jako-vytah-do-vyssich-pater'; Source = 'Moore Czech s.r.o.' }
)

$firstRow = 4
$lastRow = $firstRow + $newsItems.Count - 1
$postDateRange = "E" + $firstRow + ":E" + $lastRow

# Several PostDate strings (e.g. "02.04.2025") are day-first and would
# otherwise be auto-recognised as dates by Excels smart entry; pre-format
# the column as Text so they are kept as the literal dd.mm.yyyy strings.
$ws.Range($postDateRange).NumberFormat = "@"

$r = $firstRow
foreach ($item in $newsItems) {
    $ws.Cells.Item($r, 1).Value = $item.Title   # A: Title Original
    $ws.Cells.Item($r, 2).Value = ""            # B: Content Original (blank)
    $ws.Cells.Item($r, 3).Value = ""            # C: Title ENG (blank)
    $ws.Cells.Item($r, 4).Value = ""            # D: Summary ENG (blank)
    $ws.Cells.Item($r, 5).Value = $item.PostDate # E: PostDate
    $ws.Cells.Item($r, 6).Value = $item.Url      # F: URL
    $ws.Cells.Item($r, 7).Value = $item.Source   # G: Source
    $r++
}

